# Update table 3.4 (Stocks of Coal by Coal Rank: Electric Power Sector) from
# "... 2006 - October 2016" to "... 2006 - November 2016": refresh the title,
# add the new "November" data row before the trailing Notes row, and let the
# Notes row shift down one row as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Refresh the report title (October -> November).
$ws.Range("A1").Value2 = "Table 3.4. Stocks of Coal by Coal Rank: Electric Power Sector, 2006 - November 2016"

# 2. Insert a new row right above the merged "Notes" row (currently row 52)
#    so the Notes row shifts down to row 53, and populate it with the
#    November figures.
$ws.Rows.Item(52).Insert()

$ws.Cells.Item(52, 1).Value2 = "November"
$ws.Cells.Item(52, 2).Value2 = 71758
$ws.Cells.Item(52, 3).Value2 = 96098
$ws.Cells.Item(52, 4).Value2 = 4283
$ws.Cells.Item(52, 5).Value2 = 172139

# 3. Give the new row the same look as the other monthly data rows (e.g. the
#    October row right above it) rather than the blank formatting Insert()
#    left behind.
$ws.Range("A51:E51").Copy()
$ws.Range("A52:E52").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
